$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
  "https://www.facebook.com/congdongvnexpress/posts/941891187972893?ref=embed_post",
  "https://www.facebook.com/congdongvnexpress/posts/941906721304673?ref=embed_post",
  "https://www.facebook.com/congdongvnexpress/posts/937988701696475?ref=embed_post",
  "https://www.facebook.com/congdongvnexpress/posts/937773488384663?ref=embed_post",
  "https://www.facebook.com/congdongvnexpress/posts/936913308470681?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/952467190247614?ref=embed_post",
  "https://www.facebook.com/K14vn/posts/952328020261531?ref=embed_post",
  "https://www.facebook.com/toilarac/posts/857099016621381?ref=embed_post",
  "https://www.facebook.com/congdongvnexpress/posts/pfbid02zQ42Su3tPoWBBkg3ogkXoYQo9CGNxWAtsFqHtGxEmHraUqQGP57Hi4PCzfaCS325l",
  "https://www.facebook.com/chantroimoimedia/posts/949075047246994?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/948383670649465?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/948335790654253?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/948321197322379?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/948272710660561?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/947750444046121?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/947648517389647?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/947442614076904?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/947428554078310?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/947382224082943?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/947355524085613?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/946896847464814?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/946751607479338?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/946166314204534?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/946124337542065?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/946105157543983?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/945771730910659?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/945716020916230?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/945319524289213?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/945301487624350?ref=embed_post",
  "https://www.facebook.com/chantroimoimedia/posts/945268247627674?ref=embed_post"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $urls[$i]
}

$ws.Range("A31").Select()

$excel.ActiveWindow.Width = 13545
